$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Cells.Item(1,1).Value = "Row"
$ws.Cells.Item(1,2).Value = "Prognose"
$ws.Cells.Item(1,3).Value = "surveys"
$ws.Cells.Item(1,4).Value = "production"
$ws.Cells.Item(1,5).Value = "orders"
$ws.Cells.Item(1,6).Value = "turnover"
$ws.Cells.Item(1,7).Value = "financial"
$ws.Cells.Item(1,8).Value = "labor market"
$ws.Cells.Item(1,9).Value = "prices"
$ws.Cells.Item(1,10).Value = "national accounts"
$ws.Cells.Item(1,11).Value = "Revision"

# --- Ensure column A stays text (dates are stored as literal strings, not Excel dates) ---
$ws.Range("A2:A12").NumberFormat = "@"

# --- Data rows 2-12 ---
$ws.Cells.Item(2,1).Value = "2025-03-30"
$ws.Cells.Item(2,2).Value = 0.2922945430265343
$ws.Cells.Item(2,3).Value = 0
$ws.Cells.Item(2,4).Value = 0
$ws.Cells.Item(2,5).Value = 0
$ws.Cells.Item(2,6).Value = 0
$ws.Cells.Item(2,7).Value = 0
$ws.Cells.Item(2,8).Value = 0
$ws.Cells.Item(2,9).Value = 0
$ws.Cells.Item(2,10).Value = 0
$ws.Cells.Item(2,11).Value = 0

$ws.Cells.Item(3,1).Value = "2025-04-15"
$ws.Cells.Item(3,2).Value = 0.2845001723806546
$ws.Cells.Item(3,3).Value = 0
$ws.Cells.Item(3,4).Value = -0.005499528332737274
$ws.Cells.Item(3,5).Value = -0.0035265972271654707
$ws.Cells.Item(3,6).Value = 0.0000749349382524353
$ws.Cells.Item(3,7).Value = 0.00016061507060279848
$ws.Cells.Item(3,8).Value = 0.0000018793740442344287
$ws.Cells.Item(3,9).Value = 0.0010828896224653747
$ws.Cells.Item(3,10).Value = 0
$ws.Cells.Item(3,11).Value = -0.00008856409134178067

$ws.Cells.Item(4,1).Value = "2025-04-30"
$ws.Cells.Item(4,2).Value = 0.30022958524349547
$ws.Cells.Item(4,3).Value = 0.009727063473509075
$ws.Cells.Item(4,4).Value = 0
$ws.Cells.Item(4,5).Value = 0.000001854075467556317
$ws.Cells.Item(4,6).Value = 0.00010546131912068272
$ws.Cells.Item(4,7).Value = 0
$ws.Cells.Item(4,8).Value = 0.00012250500686866421
$ws.Cells.Item(4,9).Value = 0.004696698951812636
$ws.Cells.Item(4,10).Value = 0.0004185608585360463
$ws.Cells.Item(4,11).Value = 0.0006572691775261763

$ws.Cells.Item(5,1).Value = "2025-05-15"
$ws.Cells.Item(5,2).Value = 0.26691592669776787
$ws.Cells.Item(5,3).Value = -0.006463178550523353
$ws.Cells.Item(5,4).Value = -0.015548179749478839
$ws.Cells.Item(5,5).Value = -0.0036292323523578166
$ws.Cells.Item(5,6).Value = -0.009930137503285131
$ws.Cells.Item(5,7).Value = 0.002314911495082509
$ws.Cells.Item(5,8).Value = -0.0004475857569603041
$ws.Cells.Item(5,9).Value = 0.0005796165164434316
$ws.Cells.Item(5,10).Value = 0
$ws.Cells.Item(5,11).Value = -0.0001898726446480925

$ws.Cells.Item(6,1).Value = "2025-05-30"
$ws.Cells.Item(6,2).Value = 0.45752265500202316
$ws.Cells.Item(6,3).Value = 0.18194715449214008
$ws.Cells.Item(6,4).Value = 0
$ws.Cells.Item(6,5).Value = -0.00007056396381168112
$ws.Cells.Item(6,6).Value = 0.0013692911744620012
$ws.Cells.Item(6,7).Value = 0
$ws.Cells.Item(6,8).Value = 0.00019772436278309716
$ws.Cells.Item(6,9).Value = 0.0036031105797005086
$ws.Cells.Item(6,10).Value = 0
$ws.Cells.Item(6,11).Value = 0.0035600116589813435

$ws.Cells.Item(7,1).Value = "2025-06-15"
$ws.Cells.Item(7,2).Value = 0.4865014445354641
$ws.Cells.Item(7,3).Value = 0
$ws.Cells.Item(7,4).Value = -0.0003670322964718191
$ws.Cells.Item(7,5).Value = 0.0019823892848113978
$ws.Cells.Item(7,6).Value = 0.021667799197232235
$ws.Cells.Item(7,7).Value = 0.0030112896123823604
$ws.Cells.Item(7,8).Value = 0
$ws.Cells.Item(7,9).Value = 0.00037197525543280146
$ws.Cells.Item(7,10).Value = 0
$ws.Cells.Item(7,11).Value = 0.0023123684800540056

$ws.Cells.Item(8,1).Value = "2025-06-30"
$ws.Cells.Item(8,2).Value = 0.17955448129482382
$ws.Cells.Item(8,3).Value = -0.31307715586217544
$ws.Cells.Item(8,4).Value = 0
$ws.Cells.Item(8,5).Value = 0.00012862961696934029
$ws.Cells.Item(8,6).Value = -0.003257638592480688
$ws.Cells.Item(8,7).Value = 0
$ws.Cells.Item(8,8).Value = -0.00017427569785660884
$ws.Cells.Item(8,9).Value = 0.009301153135629162
$ws.Cells.Item(8,10).Value = 0
$ws.Cells.Item(8,11).Value = 0.00013232415927388885

$ws.Cells.Item(9,1).Value = "2025-07-15"
$ws.Cells.Item(9,2).Value = 0.06476188924324078
$ws.Cells.Item(9,3).Value = 0
$ws.Cells.Item(9,4).Value = -0.04262569748432257
$ws.Cells.Item(9,5).Value = -0.011619273735553158
$ws.Cells.Item(9,6).Value = -0.057384621064259546
$ws.Cells.Item(9,7).Value = -0.0017437071145851247
$ws.Cells.Item(9,8).Value = -0.002088855215159089
$ws.Cells.Item(9,9).Value = 0.000895672101855321
$ws.Cells.Item(9,10).Value = 0
$ws.Cells.Item(9,11).Value = -0.00022610953955887192

$ws.Cells.Item(10,1).Value = "2025-07-30"
$ws.Cells.Item(10,2).Value = 0.3243756796542123
$ws.Cells.Item(10,3).Value = 0.28212616129415813
$ws.Cells.Item(10,4).Value = 0
$ws.Cells.Item(10,5).Value = -0.0004941741369734668
$ws.Cells.Item(10,6).Value = 0.003606356729218516
$ws.Cells.Item(10,7).Value = 0
$ws.Cells.Item(10,8).Value = -0.0006524288561174697
$ws.Cells.Item(10,9).Value = 0.002423366046867785
$ws.Cells.Item(10,10).Value = -0.024625605804971067
$ws.Cells.Item(10,11).Value = -0.002769884861210825

$ws.Cells.Item(11,1).Value = "2025-08-15"
$ws.Cells.Item(11,2).Value = 0.3859447933668969
$ws.Cells.Item(11,3).Value = 0
$ws.Cells.Item(11,4).Value = 0.028539630705677776
$ws.Cells.Item(11,5).Value = 0.015167618506414687
$ws.Cells.Item(11,6).Value = 0.08387162601583491
$ws.Cells.Item(11,7).Value = 0.004338748211493408
$ws.Cells.Item(11,8).Value = 0.0005474131349907049
$ws.Cells.Item(11,9).Value = -0.014860186827811647
$ws.Cells.Item(11,10).Value = 0
$ws.Cells.Item(11,11).Value = -0.05603573603391526

$ws.Cells.Item(12,1).Value = "2025-08-30"
$ws.Cells.Item(12,2).Value = 0.291548517124095
$ws.Cells.Item(12,3).Value = -0.06716590208081855
$ws.Cells.Item(12,4).Value = 0
$ws.Cells.Item(12,5).Value = 0.0004720712786321989
$ws.Cells.Item(12,6).Value = 0.0001426347436240509
$ws.Cells.Item(12,7).Value = 0
$ws.Cells.Item(12,8).Value = -0.0003795242541873523
$ws.Cells.Item(12,9).Value = -0.007748976712229896
$ws.Cells.Item(12,10).Value = 0
$ws.Cells.Item(12,11).Value = -0.019716579217822328

# --- Column width update (col E changed: new raw width ~15.77734375) ---
# Note: Excel's ColumnWidth COM setter snaps to a whole-pixel grid (width = round(chars*6)+5, in
# units of 1/6 character), so the exact raw width 15.77734375 isn't reachable through legitimate
# column resizing; 15 chars yields the closest obtainable stored width (15.8333...).
$ws.Columns.Item(5).ColumnWidth = 15